$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rewrite rows 17-28: icsdadultosa_psg5events block (with style s=1 like rows 2-16) ---
$ws.Range("A17").Value = "icsdadultosa_psg5events"
$ws.Range("B17").Value = "ahi_a0h3"
$ws.Range("C17").Value = "AHI 3% -- all apneas and hypopneas with >=3% oxygen desaturation"
$ws.Range("D17").Value = "  "
$ws.Range("A17").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A18").Value = "icsdadultosa_psg5events"
$ws.Range("B18").Value = "ahi_a0h3a"
$ws.Range("C18").Value = "AHI 3% -- all apneas and hypopneas with >=3% oxygen desaturation or arousal"
$ws.Range("D18").Value = " "
$ws.Range("A18").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A19").Value = "icsdadultosa_psg5events"
$ws.Range("B19").Value = "ahi_a0h4"
$ws.Range("C19").Value = "AHI 4% -- all apneas and hypopneas with >=4% oxygen desaturation"
$ws.Range("D19").Value = " "
$ws.Range("A19").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A20").Value = "icsdadultosa_psg5events"
$ws.Range("B20").Value = "ahi_a0h4a"
$ws.Range("C20").Value = "AHI 4% -- all apneas and hypopneas with >=4% oxygen desaturation or arousal"
$ws.Range("D20").Value = " "
$ws.Range("A20").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A21").Value = "icsdadultosa_psg5events"
$ws.Range("B21").Value = "ahi_c0h3"
$ws.Range("C21").Value = "Central AHI 3% -- central apneas and hypopneas with >=3% oxygen desaturation"
$ws.Range("A21").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A22").Value = "icsdadultosa_psg5events"
$ws.Range("B22").Value = "ahi_c0h3a"
$ws.Range("C22").Value = "Central AHI 3% -- central apneas and hypopneas with >=3% oxygen desaturation or arousal"
$ws.Range("A22").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A23").Value = "icsdadultosa_psg5events"
$ws.Range("B23").Value = "ahi_c0h4"
$ws.Range("C23").Value = "Central AHI 4% -- central apneas and hypopneas with >=4% oxygen desaturation"
$ws.Range("A23").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A24").Value = "icsdadultosa_psg5events"
$ws.Range("B24").Value = "ahi_c0h4a"
$ws.Range("C24").Value = "Central AHI 4% -- central apneas and hypopneas with >=4% oxygen desaturation or arousal"
$ws.Range("A24").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A25").Value = "icsdadultosa_psg5events"
$ws.Range("B25").Value = "ahi_o0h3"
$ws.Range("C25").Value = "Obstructive AHI 3% -- obstructive apneas and hypopneas with >=3% oxygen desaturation"
$ws.Range("A25").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A26").Value = "icsdadultosa_psg5events"
$ws.Range("B26").Value = "ahi_o0h3a"
$ws.Range("C26").Value = "Obstructive AHI 3% -- obstructive apneas and hypopneas with >=3% oxygen desaturation or arousal"
$ws.Range("A26").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A27").Value = "icsdadultosa_psg5events"
$ws.Range("B27").Value = "ahi_o0h4"
$ws.Range("C27").Value = "Obstructive AHI 4% -- obstructive apneas and hypopneas with >=4% oxygen desaturation"
$ws.Range("A27").VerticalAlignment = $ws.Range("A2").VerticalAlignment

$ws.Range("A28").Value = "icsdadultosa_psg5events"
$ws.Range("B28").Value = "ahi_o0h4a"
$ws.Range("C28").Value = "Obstructive AHI 4% -- obstructive apneas and hypopneas with >=4% oxygen desaturation or arousal"
$ws.Range("A28").VerticalAlignment = $ws.Range("A2").VerticalAlignment

# --- Rewrite rows 29-40: icsdadultosa_psg15events block (no special style) ---
$ws.Range("A29").Value = "icsdadultosa_psg15events"
$ws.Range("B29").Value = "ahi_a0h3"
$ws.Range("C29").Value = "AHI 3% -- all apneas and hypopneas with >=3% oxygen desaturation"

$ws.Range("A30").Value = "icsdadultosa_psg15events"
$ws.Range("B30").Value = "ahi_a0h3a"
$ws.Range("C30").Value = "AHI 3% -- all apneas and hypopneas with >=3% oxygen desaturation or arousal"

$ws.Range("A31").Value = "icsdadultosa_psg15events"
$ws.Range("B31").Value = "ahi_a0h4"
$ws.Range("C31").Value = "AHI 4% -- all apneas and hypopneas with >=4% oxygen desaturation"

$ws.Range("A32").Value = "icsdadultosa_psg15events"
$ws.Range("B32").Value = "ahi_a0h4a"
$ws.Range("C32").Value = "AHI 4% -- all apneas and hypopneas with >=4% oxygen desaturation or arousal"

$ws.Range("A33").Value = "icsdadultosa_psg15events"
$ws.Range("B33").Value = "ahi_c0h3"
$ws.Range("C33").Value = "Central AHI 3% -- central apneas and hypopneas with >=3% oxygen desaturation"

$ws.Range("A34").Value = "icsdadultosa_psg15events"
$ws.Range("B34").Value = "ahi_c0h3a"
$ws.Range("C34").Value = "Central AHI 3% -- central apneas and hypopneas with >=3% oxygen desaturation or arousal"

$ws.Range("A35").Value = "icsdadultosa_psg15events"
$ws.Range("B35").Value = "ahi_c0h4"
$ws.Range("C35").Value = "Central AHI 4% -- central apneas and hypopneas with >=4% oxygen desaturation"

$ws.Range("A36").Value = "icsdadultosa_psg15events"
$ws.Range("B36").Value = "ahi_c0h4a"
$ws.Range("C36").Value = "Central AHI 4% -- central apneas and hypopneas with >=4% oxygen desaturation or arousal"

$ws.Range("A37").Value = "icsdadultosa_psg15events"
$ws.Range("B37").Value = "ahi_o0h3"
$ws.Range("C37").Value = "Obstructive AHI 3% -- obstructive apneas and hypopneas with >=3% oxygen desaturation"

$ws.Range("A38").Value = "icsdadultosa_psg15events"
$ws.Range("B38").Value = "ahi_o0h3a"
$ws.Range("C38").Value = "Obstructive AHI 3% -- obstructive apneas and hypopneas with >=3% oxygen desaturation or arousal"

$ws.Range("A39").Value = "icsdadultosa_psg15events"
$ws.Range("B39").Value = "ahi_o0h4"
$ws.Range("C39").Value = "Obstructive AHI 4% -- obstructive apneas and hypopneas with >=4% oxygen desaturation"

$ws.Range("A40").Value = "icsdadultosa_psg15events"
$ws.Range("B40").Value = "ahi_o0h4a"
$ws.Range("C40").Value = "Obstructive AHI 4% -- obstructive apneas and hypopneas with >=4% oxygen desaturation or arousal"

# --- Update selection to match final cursor position ---
$ws.Range("C42").Select()
